# PETX_YR_FIN.xlsx update — "Doing Updates for Financials"
#
# The edit inserts one new (most-recent) annual period as a new column D
# on the "PETX" worksheet, shifting the previously existing columns D:K
# one column to the right (to E:L). The three header rows (7, 38, 80)
# get a new "Period Ending" date (43465) in the new column D, and every
# other populated row gets its corresponding new-year figure in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; this shifts existing D:K -> E:L and
#    automatically extends row "spans" / sheet dimension to column L.
$ws.Columns("D").Insert()

# 2) Copy the formatting (number format / font / alignment) from the
#    column that used to be D (now E, since it kept its original look)
#    into the freshly inserted, blank column D so the new cells render
#    exactly like their neighbours (date format on header rows, number
#    format with right alignment on data rows). This is done per
#    contiguous block of populated rows so we don't manufacture empty
#    <row> elements for the blank separator rows (36 and 78) that have
#    no cells at all in either the source or the destination.
$blocks = @(
    @(7, 35),
    @(38, 77),
    @(80, 102)
)
foreach ($block in $blocks) {
    $r1 = $block[0]
    $r2 = $block[1]
    $srcFormats = $ws.Range("E${r1}:E${r2}")
    $dstFormats = $ws.Range("D${r1}:D${r2}")
    $srcFormats.Copy()
    $dstFormats.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# 3) Populate the new column D with the new period's values.
#    Key = "NA" means the text NA (shared string), matching how the
#    sheet already represents not-available figures elsewhere.
$newValues = @{
    7  = 43465
    8  = 35400
    9  = 10600
    10 = 24800
    12 = 7400
    13 = 0
    14 = 0
    15 = 500
    17 = 47300
    18 = -11900
    20 = 600
    21 = -10300
    22 = 3400
    23 = -14700
    24 = "NA"
    25 = 0
    26 = -14700
    27 = -14700
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -600
    33 = -14700
    34 = 0
    35 = -14700
    38 = 43465
    41 = 41400
    42 = 1200
    43 = 2200
    44 = 11400
    45 = 1800
    46 = 58100
    47 = "NA"
    48 = 700
    49 = 46900
    50 = 0
    51 = 0
    52 = 700
    53 = 0
    54 = 106400
    57 = 900
    58 = 0
    59 = 4600
    60 = 5600
    61 = 0
    62 = 100
    63 = 0
    64 = 0
    65 = 0
    66 = 5600
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -241200
    73 = 0
    74 = 0
    75 = 0
    76 = 100800
    77 = 0
    80 = 43465
    81 = -14700
    83 = 1000
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -11300
    91 = 0
    92 = 0
    93 = 0
    94 = -500
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = -13700
    101 = 0
    102 = -25400
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}

$wb.Save()
